$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 18 (weekly data update), shifting existing rows 18-36 down to 20-38
$ws.Rows("18:19").Insert()

# Populate new row 18
$ws.Range("A18").Value2 = 8
$ws.Range("B18").Value2 = 'Terminal La Palmera de La Serena'
$ws.Range("C18").Value2 = 'Coquimbo'
$ws.Range("D18").Value2 = 44987
$ws.Range("E18").Value2 = 4
$ws.Range("F18").Value2 = 'Fruta'
$ws.Range("G18").Value2 = 100107
$ws.Range("H18").Value2 = 'Otros'
$ws.Range("I18").Value2 = 100107011
$ws.Range("J18").Value2 = 'Tuna'
$ws.Range("K18").Value2 = 'Sin especificar'
$ws.Range("L18").Value2 = 'Especial'
$ws.Range("M18").Value2 = 160
$ws.Range("N18").Value2 = 15000
$ws.Range("O18").Value2 = 16000
$ws.Range("P18").Value2 = 15500
$ws.Range("Q18").Value2 = '$/caja 18 kilos'
$ws.Range("R18").Value2 = 'Provincia de Limarí'
$ws.Range("S18").Value2 = 861
$ws.Range("T18").Value2 = 18

# Populate new row 19
$ws.Range("A19").Value2 = 8
$ws.Range("B19").Value2 = 'Terminal La Palmera de La Serena'
$ws.Range("C19").Value2 = 'Coquimbo'
$ws.Range("D19").Value2 = 44987
$ws.Range("E19").Value2 = 4
$ws.Range("F19").Value2 = 'Fruta'
$ws.Range("G19").Value2 = 100107
$ws.Range("H19").Value2 = 'Otros'
$ws.Range("I19").Value2 = 100107011
$ws.Range("J19").Value2 = 'Tuna'
$ws.Range("K19").Value2 = 'Sin especificar'
$ws.Range("L19").Value2 = 'Primera'
$ws.Range("M19").Value2 = 240
$ws.Range("N19").Value2 = 12000
$ws.Range("O19").Value2 = 13000
$ws.Range("P19").Value2 = 12500
$ws.Range("Q19").Value2 = '$/caja 18 kilos'
$ws.Range("R19").Value2 = 'Provincia de Limarí'
$ws.Range("S19").Value2 = 694
$ws.Range("T19").Value2 = 18
